function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws 'D2' '42.332.06'
Set-TextValue $ws 'E2' '  -1.45%  '
# Row 3
Set-TextValue $ws 'D3' '2.518.63'
Set-TextValue $ws 'E3' '  -2.05%  '
# Row 5
Set-TextValue $ws 'D5' '304.26'
Set-TextValue $ws 'E5' '  +0.59%  '
# Row 6
Set-TextValue $ws 'D6' '96.88'
Set-TextValue $ws 'E6' '  -0.13%  '
# Row 7
Set-TextValue $ws 'D7' '0.587'
Set-TextValue $ws 'E7' '  +1.91%  '
# Row 8
Set-TextValue $ws 'E8' '  +0.11%  '
# Row 9
Set-TextValue $ws 'E9' '  -2.21%  '
# Row 10
Set-TextValue $ws 'D10' '36.13'
Set-TextValue $ws 'E10' '  -0.80%  '
# Row 12
Set-TextValue $ws 'D12' '0.112'
Set-TextValue $ws 'E12' '  -1.67%  '
# Row 13
Set-TextValue $ws 'E13' '  -2.20%  '
# Row 14
Set-TextValue $ws 'D14' '2.908.35'
# Row 15
Set-TextValue $ws 'B15' 'WrappedEther'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D15' '2.561.35'
Set-TextValue $ws 'E15' '  -0.70%  '
# Row 16
Set-TextValue $ws 'B16' 'Chainlink'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D16' '15.08'
Set-TextValue $ws 'E16' '  +4.83%  '
# Row 17
Set-TextValue $ws 'D17' '0.858'
Set-TextValue $ws 'E17' '  -3.01%  '
# Row 18
Set-TextValue $ws 'D18' '42.308.16'
Set-TextValue $ws 'E18' '  -1.51%  '
# Row 19
Set-TextValue $ws 'D19' '12.88'
Set-TextValue $ws 'E19' '  -1.81%  '
# Row 20
Set-TextValue $ws 'D20' '0.0₃0976'
Set-TextValue $ws 'E20' '  -1.59%  '
# Row 21
Set-TextValue $ws 'E21' '  -3.85%  '
# Row 22
Set-TextValue $ws 'D22' '71.03'
Set-TextValue $ws 'E22' '  -1.50%  '
# Row 23
Set-TextValue $ws 'D23' '251.11'
Set-TextValue $ws 'E23' '  -1.45%  '
# Row 24
Set-TextValue $ws 'D24' '2.90'
Set-TextValue $ws 'E24' '  -2.15%  '
# Row 25
Set-TextValue $ws 'D25' '2.02'
Set-TextValue $ws 'E25' '  -5.00%  '
# Row 26
Set-TextValue $ws 'E26' '  -6.50%  '
# Row 27
Set-TextValue $ws 'E27' '  +0.34%  '
# Row 28
Set-TextValue $ws 'E28' '  +9.69%  '
# Row 29
Set-TextValue $ws 'D29' '10.24'
Set-TextValue $ws 'E29' '  +0.09%  '
# Row 30
Set-TextValue $ws 'D30' '37.80'
Set-TextValue $ws 'E30' '  -0.68%  '
# Row 31
Set-TextValue $ws 'D31' '5.94'
Set-TextValue $ws 'E31' '  -2.49%  '
# Row 32
Set-TextValue $ws 'D32' '155.09'
Set-TextValue $ws 'E32' '  -0.39%  '
# Row 33
Set-TextValue $ws 'D33' '3.32'
Set-TextValue $ws 'E33' '  -1.64%  '
# Row 34
Set-TextValue $ws 'E34' '  -3.40%  '
# Row 35
Set-TextValue $ws 'E35' '  -5.27%  '
# Row 36
Set-TextValue $ws 'E36' '  -5.16%  '
# Row 37
Set-TextValue $ws 'E37' '  +0.70%  '
# Row 38
Set-TextValue $ws 'E38' '  +0.97%  '
# Row 39
Set-TextValue $ws 'E39' '  -0.28%  '
# Row 40
Set-TextValue $ws 'D40' '24.09'
Set-TextValue $ws 'E40' '  +1.57%  '
# Row 41
Set-TextValue $ws 'E41' '  -1.16%  '
# Row 42
Set-TextValue $ws 'E42' '  -1.04%  '
# Row 43
Set-TextValue $ws 'E43' '  -2.39%  '
# Row 44
Set-TextValue $ws 'D44' '0.999'
Set-TextValue $ws 'E44' '  -0.02%  '
# Row 45
Set-TextValue $ws 'B45' 'Maker'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D45' '2.045.86'
Set-TextValue $ws 'E45' '  -1.07%  '
# Row 46
Set-TextValue $ws 'B46' 'VeChain'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D46' '0.0300'
Set-TextValue $ws 'E46' '  -3.54%  '
# Row 47
Set-TextValue $ws 'D47' '84.35'
Set-TextValue $ws 'E47' '  -1.14%  '
# Row 48
Set-TextValue $ws 'D48' '8.89'
Set-TextValue $ws 'E48' '  -4.31%  '
# Row 49
Set-TextValue $ws 'D49' '2.767.42'
Set-TextValue $ws 'E49' '  -1.95%  '
# Row 50
Set-TextValue $ws 'E50' '  -1.31%  '
# Row 51
Set-TextValue $ws 'D51' '101.19'
Set-TextValue $ws 'E51' '  -4.64%  '

Write-Output "Done applying changes"
